# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the 006ef815-*.md row's
# Handoff Xliff for zh-cn was just (re-)generated, so its priority flips
# from "low" to "ht" (matching the other already-handed-off rows) and the
# handoff timestamps advance on both the zh-cn and de-de worksheets as well
# as the rolled-up "Latest HO Xliff Generate Date" on the Overview sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: row 4 is the 006ef815-8d7e-446b-8fa4-878927694f7c.md file;
# column G is "Latest HO Xliff Generate Date". It shares its value with
# rows 5-7 (89a07f4a / fe20e06e's earlier rows were already unaffected —
# those three rows track the same generation run), so all four move
# together to the new timestamp.
$wsOverview.Range("G4:G7").Value = "2016-08-29 10:30:19"

# zh-cn / de-de detail sheets: rows 4-7 correspond to the same four files.
# Priority moves from "low" to "ht" for all of them.
$wsZhCn.Range("E4:E7").Value = "ht"
$wsDeDe.Range("E4:E7").Value = "ht"

# Latest Handoff Datetime advances too - zh-cn gets a fresh timestamp,
# de-de's matches the Overview rollup value above.
$wsZhCn.Range("H4:H7").Value = "2016-08-29 10:30:00"
$wsDeDe.Range("H4:H7").Value = "2016-08-29 10:30:19"
